$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header "time_taken", matching style of the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill F2:F56 with the time_taken timestamp values (plain text, default style)
$ws.Range("F2").Value = "2021-10-05 10:52:12.448533"
$ws.Range("F3").Value = "2021-10-05 10:52:12.448546"
$ws.Range("F4").Value = "2021-10-05 10:52:12.448550"
$ws.Range("F5").Value = "2021-10-05 10:52:12.448553"
$ws.Range("F6").Value = "2021-10-05 10:52:12.448557"
$ws.Range("F7").Value = "2021-10-05 10:52:12.448560"
$ws.Range("F8").Value = "2021-10-05 10:52:12.448563"
$ws.Range("F9").Value = "2021-10-05 10:52:12.448566"
$ws.Range("F10").Value = "2021-10-05 10:52:12.448570"
$ws.Range("F11").Value = "2021-10-05 10:52:12.448573"
$ws.Range("F12").Value = "2021-10-05 10:52:12.448576"
$ws.Range("F13").Value = "2021-10-05 10:52:12.448579"
$ws.Range("F14").Value = "2021-10-05 10:52:12.448582"
$ws.Range("F15").Value = "2021-10-05 10:52:12.448585"
$ws.Range("F16").Value = "2021-10-05 10:52:12.448588"
$ws.Range("F17").Value = "2021-10-05 10:52:12.448591"
$ws.Range("F18").Value = "2021-10-05 10:52:12.448595"
$ws.Range("F19").Value = "2021-10-05 10:52:12.448598"
$ws.Range("F20").Value = "2021-10-05 10:52:12.448601"
$ws.Range("F21").Value = "2021-10-05 10:52:12.448604"
$ws.Range("F22").Value = "2021-10-05 10:52:12.448607"
$ws.Range("F23").Value = "2021-10-05 10:52:12.448610"
$ws.Range("F24").Value = "2021-10-05 10:52:12.448614"
$ws.Range("F25").Value = "2021-10-05 10:52:12.448617"
$ws.Range("F26").Value = "2021-10-05 10:52:12.448620"
$ws.Range("F27").Value = "2021-10-05 10:52:12.448624"
$ws.Range("F28").Value = "2021-10-05 10:52:12.448627"
$ws.Range("F29").Value = "2021-10-05 10:52:12.448630"
$ws.Range("F30").Value = "2021-10-05 10:52:12.448633"
$ws.Range("F31").Value = "2021-10-05 10:52:12.448636"
$ws.Range("F32").Value = "2021-10-05 10:52:12.448639"
$ws.Range("F33").Value = "2021-10-05 10:52:12.448642"
$ws.Range("F34").Value = "2021-10-05 10:52:12.448646"
$ws.Range("F35").Value = "2021-10-05 10:52:12.448649"
$ws.Range("F36").Value = "2021-10-05 10:52:12.448653"
$ws.Range("F37").Value = "2021-10-05 10:52:12.448656"
$ws.Range("F38").Value = "2021-10-05 10:52:12.448659"
$ws.Range("F39").Value = "2021-10-05 10:52:12.448662"
$ws.Range("F40").Value = "2021-10-05 10:52:12.448665"
$ws.Range("F41").Value = "2021-10-05 10:52:12.448668"
$ws.Range("F42").Value = "2021-10-05 10:52:12.448672"
$ws.Range("F43").Value = "2021-10-05 10:52:12.448675"
$ws.Range("F44").Value = "2021-10-05 10:52:12.448678"
$ws.Range("F45").Value = "2021-10-05 10:52:12.448681"
$ws.Range("F46").Value = "2021-10-05 10:52:12.448684"
$ws.Range("F47").Value = "2021-10-05 10:52:12.448687"
$ws.Range("F48").Value = "2021-10-05 10:52:12.448690"
$ws.Range("F49").Value = "2021-10-05 10:52:12.448693"
$ws.Range("F50").Value = "2021-10-05 10:52:12.448696"
$ws.Range("F51").Value = "2021-10-05 10:52:12.448699"
$ws.Range("F52").Value = "2021-10-05 10:52:12.448702"
$ws.Range("F53").Value = "2021-10-05 10:52:12.448705"
$ws.Range("F54").Value = "2021-10-05 10:52:12.448709"
$ws.Range("F55").Value = "2021-10-05 10:52:12.448712"
$ws.Range("F56").Value = "2021-10-05 10:52:12.448741"

$excel.CutCopyMode = 0
